$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price observation was recorded, so insert a new row at 65.
# This pushes the previous rows 65-130 down to 66-131 (matching the diff,
# which shows every existing row's data shifted down by one position) and
# grows the used range from A1:R130 to A1:R131.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(65, 1).Value = 3
$ws.Cells.Item(65, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(65, 3).Value = "Coquimbo"
$ws.Cells.Item(65, 4).Value = 44587
$ws.Cells.Item(65, 5).Value = 5
$ws.Cells.Item(65, 6).Value = 100112030
$ws.Cells.Item(65, 7).Value = "Poroto granado"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 65
$ws.Cells.Item(65, 11).Value = 24000
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = 24462
$ws.Cells.Item(65, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(65, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(65, 16).Value = 978
$ws.Cells.Item(65, 17).Value = 25
$ws.Cells.Item(65, 18).Value = "Hortaliza"
